$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 4737
$ws1.Range("F5").Value = 204
$ws1.Range("F8").Value = 788
$ws1.Range("F12").Value = 1161
$ws1.Range("F16").Value = 1921
$ws1.Range("F17").Value = 593
$ws1.Range("F21").Value = 221
$ws1.Range("F22").Value = 68
$ws1.Range("F23").Value = 1564
$ws1.Range("F26").Value = 2536
$ws1.Range("F30").Value = 1627
$ws1.Range("F35").Value = 4339

# Sheet "演出" (Performances) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 372
$ws2.Range("F14").Value = 51
$ws2.Range("F29").Value = 86
$ws2.Range("F38").Value = 41

# Sheet "本地生活" (Local Life) - column F update
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F7").Value = 375

# Sheet "全部类型" (All Types) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 375
$ws4.Range("F10").Value = 4737
$ws4.Range("F12").Value = 204
$ws4.Range("F15").Value = 788
$ws4.Range("F16").Value = 372
$ws4.Range("F20").Value = 1161
$ws4.Range("F26").Value = 1921
$ws4.Range("F27").Value = 593
$ws4.Range("F30").Value = 221
$ws4.Range("F31").Value = 51
$ws4.Range("F32").Value = 68
$ws4.Range("F36").Value = 1564
$ws4.Range("F41").Value = 2536
$ws4.Range("F44").Value = 1627
$ws4.Range("F49").Value = 4339
$ws4.Range("F50").Value = 41
